# "Colocando header nos gráficos"
# Adds a header label to column A (row 1) on each data sheet, drops the
# bold/border/center style from the row-label cells below the header,
# fixes the accented Portuguese spelling of several labels, removes the
# now-unused "Teto" row on the emissions sheet, and refreshes the cost
# sheet's header/values.

$wb = $excel.ActiveWorkbook

# xlPasteFormats - used to copy just the formatting (bold/border/center)
# of the existing header row onto the new A1 header cell, so it reuses
# the same style record instead of creating a new (duplicate) one.
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Sheets 1-4 ("Potencia Acumulada", "Geracao Periodo Medio",
# "Atendimento a Ponta", "Potencia Incremental") all share the same
# row layout: a header row (years) and rows 2-12 with energy-source
# labels in column A.
# ---------------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # New header for column A, formatted like the rest of row 1.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial($xlPasteFormats)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Row labels lose the bold/border/center style...
    $ws.Range("A2:A12").ClearFormats()

    # ...and a handful get corrected accents.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."
}

# ---------------------------------------------------------------------
# Sheet 5 ("Emissoes Totais"): header + P.Médio / P.Crítico rows, and
# the now-removed "Teto" row.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial($xlPasteFormats)
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2:A3").ClearFormats()
$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"

# Drop the obsolete "Teto" row entirely (was row 4).
$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6 ("Custo Total"): new header row + relabeled/updated rows.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial($xlPasteFormats)
$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Value = "2015"

$ws6.Range("A2:A3").ClearFormats()
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 573
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99

Write-Host "Header columns added and labels updated."
